$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.Value = "'97.013.97"
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '

$cell = $ws.Range('D3')
$cell.Value = "'3.672.91"
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +2.71%  '

$ws.Range('E4').Value = '  -0.08%  '

$cell = $ws.Range('D5')
$cell.Value = "'239.77"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.94%  '

$cell = $ws.Range('D6')
$cell.Value = "'1.89"
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +10.43%  '

$cell = $ws.Range('D7')
$cell.Value = "'655.63"
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.27%  '

$ws.Range('E8').Value = '  -0.77%  '

$ws.Range('E9').Value = '  +3.92%  '

$ws.Range('E10').Value = '  +0.02%  '

$cell = $ws.Range('D11')
$cell.Value = "'3.671.94"
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +2.79%  '

$cell = $ws.Range('D12')
$cell.Value = "'45.54"
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +2.67%  '

$ws.Range('E13').Value = '  +1.06%  '

$cell = $ws.Range('D14')
$cell.Value = "'6.85"
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +6.59%  '

$cell = $ws.Range('D15')
$cell.Value = "'4.357.10"
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +2.72%  '

$ws.Range('E16').Value = '  +3.38%  '

$cell = $ws.Range('D17')
$cell.Value = "'96.681.22"
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -0.32%  '

$cell = $ws.Range('D18')
$cell.Value = "'8.97"
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +3.80%  '

$cell = $ws.Range('D19')
$cell.Value = "'3.687.04"
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +3.30%  '

$cell = $ws.Range('D20')
$cell.Value = "'18.95"
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +5.34%  '

$cell = $ws.Range('D21')
$cell.Value = "'12.79"
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.66%  '

$cell = $ws.Range('D22')
$cell.Value = "'0.533"
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +1.34%  '

$cell = $ws.Range('D23')
$cell.Value = "'532.79"
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +3.44%  '

$ws.Range('E24').Value = '  +0.34%  '

$cell = $ws.Range('D25')
$cell.Value = "'7.21"
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +5.50%  '

$ws.Range('E26').Value = '  -0.40%  '

$cell = $ws.Range('D27')
$cell.Value = "'102.47"
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +1.14%  '

$cell = $ws.Range('D28')
$cell.Value = "'13.51"
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +3.63%  '

$ws.Range('E29').Value = '  +1.87%  '

$cell = $ws.Range('D30')
$cell.Value = "'12.65"
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +6.61%  '

$ws.Range('E31').Value = '  +2.15%  '

$cell = $ws.Range('D32')
$cell.Value = "'0.999"
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('E33').Value = '  +14.65%  '

$ws.Range('E34').Value = '  +1.05%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D35')
$cell.Value = "'32.75"
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +3.41%  '

$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell = $ws.Range('D36')
$cell.Value = "'1.00"
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +0.24%  '

$cell = $ws.Range('D37')
$cell.Value = "'653.65"
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +5.73%  '

$ws.Range('E38').Value = '  +5.62%  '

$cell = $ws.Range('D39')
$cell.Value = "'8.91"
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.51%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D40')
$cell.Value = "'6.89"
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +14.97%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range('D41')
$cell.Value = "'0.163"
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +5.66%  '

$ws.Range('E42').Value = '  +2.94%  '

$ws.Range('E43').Value = '  +4.47%  '

$cell = $ws.Range('D44')
$cell.Value = "'38.80"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +18.03%  '

$ws.Range('E45').Value = '  +0.02%  '

$cell = $ws.Range('D46')
$cell.Value = "'0.0460"
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +4.87%  '

$cell = $ws.Range('D47')
$cell.Value = "'0.443"
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +9.22%  '

$ws.Range('E48').Value = '  +1.73%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range('D49')
$cell.Value = "'8.78"
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +3.36%  '

$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$cell = $ws.Range('D50')
$cell.Value = "'23.63"
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +0.16%  '

$cell = $ws.Range('D51')
$cell.Value = "'3.64"
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +4.08%  '
